# Adding support for Parallel Cross Browser Testing
#
# The DATA worksheet's test-data table gains a new "browser" column
# (populated with "firefox" for every existing data row) inserted right
# after the "execute" column and before "username" - i.e. at column C,
# pushing username/password/firstname one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Insert a new blank column at C; existing C:E (username/password/firstname)
# shift right to D:F, exactly mirroring Excel's native "Insert" column
# command (dimension, column widths and styles all shift automatically).
$ws.Columns("C").Insert()

# Populate the new "browser" column.
$ws.Range("C1").Value = "browser"
$ws.Range("C2").Value = "firefox"
$ws.Range("C3").Value = "firefox"
$ws.Range("C4").Value = "firefox"
$ws.Range("C5").Value = "firefox"
$ws.Range("C6").Value = "firefox"
